# Table-X-SurfaceVariables.xlsx
# "Added temp and fluor mean profiles for STFZ. Changed station list"
#
# The "Stations and Times" sheet currently holds a single By-EYE station
# table (rows 2-5). This adds a section label above it, duplicates the
# table for a new "By PROFILE" data set (rows 9-12), and updates the
# station/latitude columns of the new table with the STFZ profile values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stations and Times")

# Section label above the existing ("by eye") table.
$ws.Range("A1").Value = "By EYE"

# Section label above the new ("by profile") table, mirroring row 2's header.
$ws.Range("A8").Value = "By PROFILE"

# Duplicate the data rows (3:5) -> (10:12), carrying over the same
# formatting (date styles, borders, etc.) used by the first table.
$ws.Range("A3:F5").Copy($ws.Range("A10"))

# New station / latitude values for the duplicated ("by profile") rows;
# the dates/distances/times columns (C:F) stay identical to rows 3:5.
$ws.Range("A10").Value = "27-29"
$ws.Range("B10").Value = "32.25°-32.75°N"

$ws.Range("A11").Value = "18-22"
$ws.Range("B11").Value = "31.25°-32.25°N"

$ws.Range("A12").Value = "21-25"
$ws.Range("B12").Value = "31.25°-32.25°N"

$ws.Range("B12").Select()
